$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data gained one more entry (2026/02/03) ahead of the existing
# 2026/12/29.. block, so every row from 747 through 788 shifts down by one.
# Inserting a whole row at 747 reproduces that shift (and keeps the sheet's
# <dimension> in sync automatically).
$ws.Rows(747).Insert()

# Text columns A/B would otherwise be auto-coerced (e.g. "2026/02/03" parsed
# as a real date serial). Force them to plain text first, write the values,
# then drop the temporary "@" format so the cells end up unstyled - same as
# every other date/weekday cell in this column.
$textCells = $ws.Range("A747:B747")
$textCells.NumberFormat = "@"
$ws.Range("A747").Value = "2026/02/03"
$ws.Range("B747").Value = "火"
$textCells.ClearFormats()

$ws.Range("C747").Value = 13
$ws.Range("D747").Value = 23
